# A new daily price record (Provincia de Curicó, 2023-01-30) was added to the
# "Ciruela" sheet just below the header, pushing every existing record down by
# one row. Reproduce this by inserting a new row 3, seeding it with a copy of
# the row that lands there after the shift (i.e. the original row 3's data),
# then overwriting the date and origin to the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 3; everything from old row 3 onward moves down
# by one (old row 3 -> new row 4, ..., old row 74 -> new row 75).
$ws.Rows.Item(3).Insert()

# Seed the new row 3 with the same Mercado/Producto/Variedad/Calidad/Volumen/
# Precio/Unidad/Kg values as the record now sitting in row 4 (the former
# row 3), then fix up the date and origin for the new record.
$ws.Range("A3:T3").Value2 = $ws.Range("A4:T4").Value2

$ws.Range("D3").Value2 = 44956
$ws.Range("R3").Value2 = "Provincia de Curicó"
